# Completed the README file.
# Fills in the developer name and the test-plan rows (Method Inputs /
# Expected Result columns) for the Client unit-test plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Nishant Malhotra"

# Test case 1 (__init__, valid attributes)
$ws.Range("F7").Value = "1010, 'Susan', 'Clark', 'susanclark@pixell.com'`t"
$ws.Range("G7").Value = "Attributes set correctly; no error raised."

# Test case 2 (__init__, invalid client number)
# Leading apostrophe marks the value as text (quote-prefix) — matches the
# source data which begins with an unmatched quote character.
$ws.Range("F8").Value = "'ABC', 'Susan', 'Clark', 'susanclark@pixell.com'`t"
$ws.Range("G8").Value = "ValueError: Client number must be an integer."

# Test case 3 (__init__, blank first name)
$ws.Range("F9").Value = "1010, ' ', 'Clark', 'susanclark@pixell.com'`t"
$ws.Range("G9").Value = "ValueError: First name cannot be blank."

# Test case 4 (__init__, blank last name)
$ws.Range("F10").Value = "1010, 'Susan', ' ', 'susanclark@pixell.com'`t"
$ws.Range("G10").Value = "ValueError: Last name cannot be blank."

# Test case 5 (__init__, invalid email)
$ws.Range("F11").Value = "1010, 'Susan', 'Clark', 'invalid-email'`t"
$ws.Range("G11").Value = "Email address is set to email@pixell-river.com."

# Test case 6 (client_number getter)
$ws.Range("F12").Value = "N/A`t"
$ws.Range("G12").Value = "Returns 1010."

# Test case 7 (first_name getter)
$ws.Range("F13").Value = "N/A`t"
$ws.Range("G13").Value = "Returns 'Susan'."

# Test case 8 (last_name getter)
$ws.Range("F14").Value = "N/A`t"
$ws.Range("G14").Value = "Returns 'Clark'."

# Test case 9 (email_address getter)
$ws.Range("F15").Value = "N/A`t"
$ws.Range("G15").Value = "Returns 'susanclark@pixell.com'."

# Test case 10 (__str__)
$ws.Range("F16").Value = "N/A`tReturns 'Clark, Susan [1010] - susanclark@pixell.com`n'"
$ws.Range("G16").Value = "Returns 'Clark, Susan [1010] - susanclark@pixell.com`n'"

# Reflect the last cell the author was working in / scrolled to.
$ws.Range("G8").Select()
